# Update cryptocurrency price/volume data per Mon May 15 18:29:32 UTC 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.688.91"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "'1.846.78"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  -1.72%  "

$ws.Range("D5").Value = "'318.16"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("D7").Value = "'0.4301"
$ws.Range("E7").Value = "  -1.64%  "

$ws.Range("D8").Value = "'0.3750"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'0.07329"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("D10").Value = "'0.8747"
$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("D11").Value = "'21.51"

$ws.Range("D12").Value = "'1.847.94"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "'6.741"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "'5.432"
$ws.Range("E14").Value = "  -1.03%  "

$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "'89.05"
$ws.Range("E16").Value = "  +4.65%  "

$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").Value = "'0.000009004"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("D20").Value = "'15.42"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'27.706.65"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").Value = "'5.208"
$ws.Range("E22").Value = "  -1.37%  "

$ws.Range("D23").Value = "'11.03"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").Value = "'2.073.53"
$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("D26").Value = "'155.39"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").Value = "'18.62"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").Value = "'2.165"
$ws.Range("E28").Value = "  +8.91%  "

$ws.Range("D29").Value = "'5.348"
$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").Value = "'118.80"
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("D31").Value = "'0.08929"
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").Value = "'1.225"
$ws.Range("E32").Value = "  +1.54%  "

$ws.Range("D33").Value = "'0.7767"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").Value = "'4.539"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("D35").Value = "'2.920"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").Value = "'1.131"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").Value = "'0.01982"

$ws.Range("D39").Value = "'0.05337"
$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("E40").Value = "  +2.05%  "

$ws.Range("D41").Value = "'7.148"
$ws.Range("E41").Value = "  +4.64%  "

$ws.Range("D42").Value = "'0.1695"
$ws.Range("E42").Value = "  +1.65%  "

$ws.Range("D43").Value = "'0.5130"
$ws.Range("E43").Value = "  -0.72%  "

$ws.Range("D44").Value = "'8.781"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'10.72"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("D46").Value = "'107.51"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").Value = "'0.4758"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("D48").Value = "'0.06459"
$ws.Range("E48").Value = "  -2.05%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.012"
$ws.Range("E49").Value = "  -1.64%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.688"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").Value = "'1.846"
$ws.Range("E51").Value = "  -2.32%  "
